$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("A18").Value = "Maplin"
$ws.Range("C18").Value = "Maplin 4 Pole 3.5mm Stereo Jack Cable 0.75m"
$ws.Range("B18").Value = "A61NW"
$ws.Range("D18").Value = 11.99
$ws.Range("E18").Value = 3
$ws.Range("F18").Formula = "=IF(NOT(ISBLANK(`$D18)),`$D18*`$E18,`"`")"

# Row 19
$ws.Range("A19").Value = "RS"
$ws.Range("C19").Value = "HARWIN D01-99 Series, 2.54mm Pitch 20 Way 1 Row Straight PCB Header, Solder Termination, 1.5A"
$ws.Range("B19").Value = "547-3302"
$ws.Range("D19").Value = 13.79
$ws.Range("E19").Value = 5
$ws.Range("F19").Formula = "=IF(NOT(ISBLANK(`$D19)),`$D19*`$E19,`"`")"

# Row 20
$ws.Range("A20").Value = "Farnell"
$ws.Range("C20").Value = "KINGSTATE  KEEG1542PBL-A  MICROPHONE, ELECTRET, CONDENSER"
$ws.Range("B20").Value = 1502746
$ws.Range("D20").Value = 1.12
$ws.Range("E20").Value = 5
$ws.Range("F20").Formula = "=IF(NOT(ISBLANK(`$D20)),`$D20*`$E20,`"`")"

# Row 21
$ws.Range("A21").Value = "RS"
$ws.Range("C21").Value = "AD605BRZ, Dual Voltage Controlled Amplifier Single Ended 4.5 → 5.5 V 16-Pin SOIC"
$ws.Range("B21").Value = "758-9705"
$ws.Range("D21").Value = 13.91
$ws.Range("E21").Value = 5
$ws.Range("F21").Formula = "=IF(NOT(ISBLANK(`$D21)),`$D21*`$E21,`"`")"

# Row 22
$ws.Range("A22").Value = "Farnell"
$ws.Range("C22").Value = "LUMBERG  1503 13 VP3  SOCKET, 3.5MM JACK, SMT, 4WAY"
$ws.Range("B22").Value = 1368640
$ws.Range("D22").Value = 0.561
$ws.Range("E22").Value = 5
$ws.Range("F22").Formula = "=IF(NOT(ISBLANK(`$D22)),`$D22*`$E22,`"`")"

# Row 23
$ws.Range("A23").Value = "Farnell"
$ws.Range("C23").Value = "VISHAY BEYSCHLAG  MMA02040E1001BB100  RES, MELF, 1K, 0.1%, 250MW, SMD"
$ws.Range("B23").Value = 3085946
$ws.Range("D23").Value = 0.214
$ws.Range("E23").Value = 5
$ws.Range("F23").Formula = "=IF(NOT(ISBLANK(`$D23)),`$D23*`$E23,`"`")"

# Row 24
$ws.Range("A24").Value = "Farnell"
$ws.Range("C24").Value = "VISHAY BEYSCHLAG  MMA02040C2201FB300  RES, MELF, 2K2, 1%, 250MW, SMD"
$ws.Range("B24").Value = 3087827
$ws.Range("D24").Value = 0.0157
$ws.Range("E24").Value = 10
$ws.Range("F24").Formula = "=IF(NOT(ISBLANK(`$D24)),`$D24*`$E24,`"`")"

# Row 25
$ws.Range("A25").Value = "Farnell"
$ws.Range("C25").Value = "VISHAY BEYSCHLAG  MMA02040E1002BB100  RESISTOR, 0204 10K"
$ws.Range("B25").Value = "3086185RL"
$ws.Range("D25").Value = 0.209
$ws.Range("E25").Value = 10
$ws.Range("F25").Formula = "=IF(NOT(ISBLANK(`$D25)),`$D25*`$E25,`"`")"

# Row 26
$ws.Range("A26").Value = "Farnell"
$ws.Range("C26").Value = "VISHAY BEYSCHLAG  MMA02040C6802FB300  RES, MELF, 68K, 1%, 250MW, SMD"
$ws.Range("B26").Value = 3088005
$ws.Range("D26").Value = 0.0178
$ws.Range("E26").Value = 5
$ws.Range("F26").Formula = "=IF(NOT(ISBLANK(`$D26)),`$D26*`$E26,`"`")"

# Row 27
$ws.Range("A27").Value = "Farnell"
$ws.Range("C27").Value = "VISHAY BEYSCHLAG  MMA02040E8202BB100  RESISTOR, 0204 82K"
$ws.Range("B27").Value = "3086409RL"
$ws.Range("D27").Value = 0.258
$ws.Range("E27").Value = 5
$ws.Range("F27").Formula = "=IF(NOT(ISBLANK(`$D27)),`$D27*`$E27,`"`")"

# Row 28
$ws.Range("A28").Value = "Farnell"
$ws.Range("C28").Value = "TDK  C3216X7R1C106K160AC  CAP, MLCC, X7R, 10UF, 16V, 1206"
$ws.Range("B28").Value = "1907353RL"
$ws.Range("D28").Value = 0.0364
$ws.Range("E28").Value = 25
$ws.Range("F28").Formula = "=IF(NOT(ISBLANK(`$D28)),`$D28*`$E28,`"`")"

# Row 29
$ws.Range("A29").Value = "Farnell"
$ws.Range("C29").Value = "TDK  C3216C0G2J101J060AA  CERAMIC CAPACITOR 100PF 630V, C0G, 5%, 1206"
$ws.Range("B29").Value = 1844418
$ws.Range("D29").Value = 0.02
$ws.Range("E29").Value = 10
$ws.Range("F29").Formula = "=IF(NOT(ISBLANK(`$D29)),`$D29*`$E29,`"`")"

# Row 30
$ws.Range("A30").Value = "Farnell"
$ws.Range("C30").Value = "CHEMTRONICS  CW8400  DISPENSING PEN, FLUX, LEAD FREE, 9G"
$ws.Range("B30").Value = 9599568
$ws.Range("D30").Value = 8.51
$ws.Range("E30").Value = 1
$ws.Range("F30").Formula = "=IF(NOT(ISBLANK(`$D30)),`$D30*`$E30,`"`")"

# Final row: mark as "PAID MYSELF" with a red highlight fill
$ws.Range("G18").Value = "PAID MYSELF"
$ws.Range("G18").Interior.Color = 255

# Column G width (new column added to the right of the table)
$ws.Columns("G").ColumnWidth = 12.33

# Recalculate all formulas (SUM(F:F) on G2 and the new F column entries)
$excel.Calculate()

# Update the view: scroll so row 7 is at the top, and select G20 (matches author's final cursor position)
$ws.Range("G20").Select()
$excel.ActiveWindow.ScrollRow = 7
